# Applies the numeric updates to the Pandaemonium_Profits workbook
# (rows across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1686.28
$ws.Range("J112").Value = 1798.1305
$ws.Range("L112").Value = 5394.3915
$ws.Range("N112").Value = -7610.3915

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1823829.5
$ws.Range("I138").Value = 4287.1
$ws.Range("J138").Value = 2278715
$ws.Range("K138").Value = 12861.3
$ws.Range("L138").Value = 6836145
$ws.Range("M138").Value = -7721.300000000001
$ws.Range("N138").Value = -6846425

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1578.6
$ws.Range("I45").Value = 1545.24
$ws.Range("J45").Value = 1745.4
$ws.Range("K45").Value = 1545.24
$ws.Range("L45").Value = 1745.4
$ws.Range("M45").Value = -1168.24
$ws.Range("N45").Value = -2499.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6673.979
$ws.Range("I61").Value = 3344.5144
$ws.Range("J61").Value = 15637.923
$ws.Range("K61").Value = 3344.5144
$ws.Range("L61").Value = 15637.923
$ws.Range("M61").Value = -3132.5144
$ws.Range("N61").Value = -16061.923

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I122").Value = 1828.7273
$ws.Range("K122").Value = 5486.1819
$ws.Range("M122").Value = -3036.1819

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4576.66
$ws.Range("I132").Value = 1904.4231
$ws.Range("J132").Value = 7149.926
$ws.Range("K132").Value = 5713.2693
$ws.Range("L132").Value = 21449.778
$ws.Range("M132").Value = -3183.2693
$ws.Range("N132").Value = -26509.778

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 6673.979
$ws.Range("I136").Value = 3344.5144
$ws.Range("J136").Value = 15637.923
$ws.Range("K136").Value = 10033.5432
$ws.Range("L136").Value = 46913.769
$ws.Range("M136").Value = -7483.5432
$ws.Range("N136").Value = -52013.769

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 20095.387
$ws.Range("I134").Value = 2222.2954
$ws.Range("K134").Value = 6666.8862
$ws.Range("M134").Value = -4131.8862

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4336.511
$ws.Range("I31").Value = 1782.5714
$ws.Range("J31").Value = 5489.9033
$ws.Range("K31").Value = 1782.5714
$ws.Range("L31").Value = 5489.9033
$ws.Range("M31").Value = -1487.5714
$ws.Range("N31").Value = -6079.9033

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4336.511
$ws.Range("I34").Value = 1782.5714
$ws.Range("J34").Value = 5489.9033
$ws.Range("K34").Value = 1782.5714
$ws.Range("L34").Value = 5489.9033
$ws.Range("M34").Value = -1580.5714
$ws.Range("N34").Value = -5893.9033

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 15452.883
$ws.Range("I122").Value = 8687.429
$ws.Range("K122").Value = 26062.287
$ws.Range("M122").Value = -23612.287

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2412.3076
$ws.Range("I132").Value = 2121.1
$ws.Range("J132").Value = 3383
$ws.Range("K132").Value = 6363.299999999999
$ws.Range("L132").Value = 10149
$ws.Range("M132").Value = -3833.299999999999
$ws.Range("N132").Value = -15209

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3731.16
$ws.Range("I134").Value = 2686.6
$ws.Range("J134").Value = 4427.533
$ws.Range("K134").Value = 8059.799999999999
$ws.Range("L134").Value = 13282.599
$ws.Range("M134").Value = -5524.799999999999
$ws.Range("N134").Value = -18352.599

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1458.7142
$ws.Range("I68").Value = 800
$ws.Range("J68").Value = 1568.5
$ws.Range("K68").Value = 2400
$ws.Range("L68").Value = 4705.5
$ws.Range("M68").Value = -1589
$ws.Range("N68").Value = -6327.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 14707211
$ws.Range("J69").Value = 16130441
$ws.Range("L69").Value = 48391323
$ws.Range("N69").Value = -48392945

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 1458.7142
$ws.Range("I71").Value = 800
$ws.Range("J71").Value = 1568.5
$ws.Range("K71").Value = 7200
$ws.Range("L71").Value = 14116.5
$ws.Range("M71").Value = -3144
$ws.Range("N71").Value = -22228.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H72").Value = 14707211
$ws.Range("J72").Value = 16130441
$ws.Range("L72").Value = 145173969
$ws.Range("N72").Value = -145182081

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 2913.5715
$ws.Range("J81").Value = 2868.4614
$ws.Range("L81").Value = 8605.3842
$ws.Range("N81").Value = -10851.3842

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H84").Value = 2913.5715
$ws.Range("J84").Value = 2868.4614
$ws.Range("L84").Value = 25816.1526
$ws.Range("N84").Value = -37048.1526

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 5146.625
$ws.Range("J133").Value = 4931.364
$ws.Range("L133").Value = 14794.092
$ws.Range("N133").Value = -24914.092

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 6200.1055
$ws.Range("I102").Value = 5889.091
$ws.Range("J102").Value = 6627.75
$ws.Range("K102").Value = 5889.091
$ws.Range("L102").Value = 6627.75
$ws.Range("M102").Value = -4267.091
$ws.Range("N102").Value = -9871.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2542.6
$ws.Range("I113").Value = 2925
$ws.Range("J113").Value = 1013
$ws.Range("K113").Value = 2925
$ws.Range("L113").Value = 1013
$ws.Range("M113").Value = -755
$ws.Range("N113").Value = -5353

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 13571.286
$ws.Range("I122").Value = 14000
$ws.Range("J122").Value = 12999.667
$ws.Range("K122").Value = 42000
$ws.Range("L122").Value = 38999.001
$ws.Range("M122").Value = -39550
$ws.Range("N122").Value = -43899.001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 6454.919
$ws.Range("I132").Value = 6277.75
$ws.Range("J132").Value = 6503.793
$ws.Range("K132").Value = 18833.25
$ws.Range("L132").Value = 19511.379
$ws.Range("M132").Value = -16303.25
$ws.Range("N132").Value = -24571.379

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 690.6842
$ws.Range("I16").Value = 562.2778
$ws.Range("J16").Value = 3002
$ws.Range("K16").Value = 562.2778
$ws.Range("L16").Value = 3002
$ws.Range("M16").Value = -392.2778
$ws.Range("N16").Value = -3342

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6008.4644
$ws.Range("I122").Value = 4661.3335
$ws.Range("J122").Value = 8433.299999999999
$ws.Range("K122").Value = 13984.0005
$ws.Range("L122").Value = 25299.9
$ws.Range("M122").Value = -11534.0005
$ws.Range("N122").Value = -30199.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3141.8076
$ws.Range("I122").Value = 1184.5
$ws.Range("J122").Value = 5425.3335
$ws.Range("K122").Value = 3553.5
$ws.Range("L122").Value = 16276.0005
$ws.Range("M122").Value = -1103.5
$ws.Range("N122").Value = -21176.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1789.55
$ws.Range("I132").Value = 1237.52
$ws.Range("J132").Value = 2709.6
$ws.Range("K132").Value = 3712.56
$ws.Range("L132").Value = 8128.799999999999
$ws.Range("M132").Value = -1182.56
$ws.Range("N132").Value = -13188.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 5325.2
$ws.Range("I136").Value = 2975.9167
$ws.Range("J136").Value = 8241.552
$ws.Range("K136").Value = 8927.750100000001
$ws.Range("L136").Value = 24724.656
$ws.Range("M136").Value = -6377.750100000001
$ws.Range("N136").Value = -29824.656
